# Refine metadata to be an additional tab.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- 1. Update the "time_taken" timestamps on the "data" sheet (column F) ---
$data.Range("F2").Value = "2021-10-05 14:22:28.371484"
$data.Range("F3").Value = "2021-10-05 14:22:28.371492"
$data.Range("F4").Value = "2021-10-05 14:22:28.371495"
$data.Range("F5").Value = "2021-10-05 14:22:28.371498"
$data.Range("F6").Value = "2021-10-05 14:22:28.371501"
$data.Range("F7").Value = "2021-10-05 14:22:28.371504"
$data.Range("F8").Value = "2021-10-05 14:22:28.371506"
$data.Range("F9").Value = "2021-10-05 14:22:28.371511"
$data.Range("F10").Value = "2021-10-05 14:22:28.371514"
$data.Range("F11").Value = "2021-10-05 14:22:28.371517"

# --- 2. Add a new "metadata" worksheet right after the "data" sheet ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row values
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row 2 values
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Renal cancer pertinent cancer susceptibility"
$meta.Range("C2").Value = 154

$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.1"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2019-06-20T15:14:12.998750Z"
$meta.Range("F2").Value = "2021-10-05 14:22:28.367823"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/154/?format=json"

# Apply the same header / "index column" cell style already used on the
# "data" sheet (bold, thin border, centered/top aligned) by copying the
# format from the existing styled cells rather than rebuilding it property
# by property (which would create a brand-new, slightly different style).
$data.Range("B1").Copy() | Out-Null
$meta.Range("B1:G1").PasteSpecial(-4122) | Out-Null

$data.Range("A2").Copy() | Out-Null
$meta.Range("A2").PasteSpecial(-4122) | Out-Null

$wb.Application.CutCopyMode = $false

$meta.Range("A1").Select() | Out-Null
